$wb = $excel.ActiveWorkbook

# --- Sheet "Data" ---
$wsData = $wb.Worksheets.Item("Data")

# Row 2
$wsData.Range("A2").Value = 3035
$wsData.Range("E2").Value = 46200608035
$wsData.Range("X2").Value = "DN4127460130032"

# Row 3
$wsData.Range("A3").Value = 3036
$wsData.Range("E3").Value = 46200608036

# Row 4
$wsData.Range("A4").Value = 3037
$wsData.Range("E4").Value = 46200608037
$wsData.Range("X4").Value = "DN4127460130034"

# Row 5
$wsData.Range("A5").Value = 3038
$wsData.Range("E5").Value = 46200608038

# --- Sheet "Check" ---
$wsCheck = $wb.Worksheets.Item("Check")

# Row 2
$wsCheck.Range("A2").Value = 3035
$wsCheck.Range("C2").Value = "DN4127460130032"

# Row 3
$wsCheck.Range("A3").Value = 3036

# Row 4
$wsCheck.Range("A4").Value = 3037
$wsCheck.Range("C4").Value = "DN4127460130034"

# Row 5
$wsCheck.Range("A5").Value = 3038
